$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.091.54'
$ws.Range('E2').Value = '  -1.83%  '
$ws.Range('D3').Value = '1.835.20'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '325.66'
$ws.Range('E5').Value = '  -3.02%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.4615'
$ws.Range('E7').Value = '  -1.04%  '
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('D9').Value = '0.07846'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').Value = '0.9606'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').Value = '21.99'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('D12').Value = '1.850.39'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '5.676'
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').Value = '6.898'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '88.34'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '0.000009933'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').Value = '16.68'
$ws.Range('E19').Value = '  -2.37%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '28.088.24'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('D23').Value = '11.00'
$ws.Range('E23').Value = '  -2.35%  '
$ws.Range('D24').Value = '2.086'
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('D25').Value = '2.041.75'
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D26').Value = '154.62'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').Value = '19.11'
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('D28').Value = '5.730'
$ws.Range('E28').Value = '  -5.74%  '
$ws.Range('D29').Value = '1.976'
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('D30').Value = '119.47'
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('D31').Value = '0.9379'
$ws.Range('E31').Value = '  -3.40%  '
$ws.Range('D32').Value = '0.09236'
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('D33').Value = '5.268'
$ws.Range('E33').Value = '  -1.73%  '
$ws.Range('D34').Value = '1.319'
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').Value = '3.332'
$ws.Range('E35').Value = '  -4.30%  '
$ws.Range('D36').Value = '0.05826'
$ws.Range('E36').Value = '  -5.55%  '
$ws.Range('D37').Value = '0.02117'
$ws.Range('E37').Value = '  -3.80%  '
$ws.Range('D38').Value = '1.137'
$ws.Range('E38').Value = '  -2.81%  '
$ws.Range('D39').Value = '7.711'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').Value = '0.5601'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('D41').Value = '9.940'
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('D42').Value = '0.1756'
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').Value = '0.07320'
$ws.Range('E43').Value = '  +3.12%  '
$ws.Range('D44').Value = '11.70'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').Value = '0.5267'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('D46').Value = '2.129'
$ws.Range('E46').Value = '  -11.60%  '
$ws.Range('E47').Value = '  -7.22%  '
$ws.Range('D48').Value = '1.836'
$ws.Range('E48').Value = '  -3.62%  '
$ws.Range('D49').Value = '113.71'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').Value = '1.019'
$ws.Range('E51').Value = '  -0.11%  '
